$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new value reads as a plain number -- pre-format as
# text so they keep being stored as strings, like the source file.
$textPriceRows = @(4,5,6,7,11,12,13,14,16,20,22,23,24,25,26,27,28,29,30,31,32,33,34,35,37,38,39,40,44,45,46,47,48,49,50,51)
foreach ($r in $textPriceRows) { $ws.Range("D$r").NumberFormat = "@" }

# --- Row-by-row price (D) / volume (E) updates ---
$ws.Range("D2").Value = '66.875.74'
$ws.Range("E2").Value = '  +2.90%  '
$ws.Range("D3").Value = '3.441.71'
$ws.Range("E3").Value = '  +1.59%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = '570.90'
$ws.Range("E5").Value = '  +1.96%  '
$ws.Range("D6").Value = '183.92'
$ws.Range("E6").Value = '  +5.70%  '
$ws.Range("D7").Value = '0.635'
$ws.Range("E7").Value = '  +1.82%  '
$ws.Range("D8").Value = '3.439.65'
$ws.Range("E8").Value = '  +1.81%  '
$ws.Range("E9").Value = '  -0.15%  '
$ws.Range("E10").Value = '  +6.68%  '
$ws.Range("D11").Value = '0.645'
$ws.Range("E11").Value = '  +2.17%  '
$ws.Range("D12").Value = '55.25'
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("D13").Value = '0.0000281'
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("D14").Value = '9.38'
$ws.Range("E14").Value = '  +3.06%  '
$ws.Range("D15").Value = '3.971.48'
$ws.Range("E15").Value = '  +1.02%  '
$ws.Range("D16").Value = '18.51'
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").Value = '3.426.15'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("E18").Value = '  +0.47%  '
$ws.Range("D19").Value = '66.700.31'
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("D20").Value = '12.06'
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").Value = '470.85'
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '5.00'
$ws.Range("E23").Value = '  +2.10%  '
$ws.Range("D24").Value = '14.74'
$ws.Range("E24").Value = '  +9.08%  '
$ws.Range("D25").Value = '4.18'
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("D26").Value = '89.94'
$ws.Range("E26").Value = '  +3.52%  '
$ws.Range("D27").Value = '2.95'
$ws.Range("D28").Value = '10.94'
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("D29").Value = '8.93'
$ws.Range("E29").Value = '  +1.98%  '
$ws.Range("D30").Value = '31.54'
$ws.Range("E30").Value = '  +2.69%  '
$ws.Range("D31").Value = '6.96'
$ws.Range("E31").Value = '  +2.85%  '
$ws.Range("D32").Value = '11.64'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("D33").Value = '588.44'
$ws.Range("E33").Value = '  +2.73%  '
$ws.Range("D34").Value = '62.63'
$ws.Range("E34").Value = '  +1.93%  '
$ws.Range("D35").Value = '0.110'
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("E36").Value = '  -0.11%  '
$ws.Range("D37").Value = '0.147'
$ws.Range("E37").Value = '  +5.92%  '
$ws.Range("D38").Value = '3.64'
$ws.Range("E38").Value = '  +2.33%  '
$ws.Range("D39").Value = '36.60'
$ws.Range("E39").Value = '  +3.01%  '
$ws.Range("D40").Value = '0.388'
$ws.Range("E40").Value = '  +5.03%  '
$ws.Range("D41").Value = '0.0₃0766'
$ws.Range("E41").Value = '  +3.18%  '
$ws.Range("D42").Value = '3.132.59'
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("E43").Value = '  +3.07%  '
$ws.Range("D44").Value = '0.0427'
$ws.Range("E44").Value = '  +2.77%  '
$ws.Range("D49").Value = '0.998'
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("D50").Value = '141.69'
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("D51").Value = '8.67'
$ws.Range("E51").Value = '  +4.84%  '

# --- Rows 45-48 reordered/replaced (coin rank shuffle) ---
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '2.85'
$ws.Range("E45").Value = '  +22.56%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '2.55'
$ws.Range("E46").Value = '  +3.75%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '3.21'
$ws.Range("E47").Value = '  +1.86%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '0.135'
$ws.Range("E48").Value = '  +0.32%  '

# Restore default (Normal) style on the price cells we forced to text,
# so no stray number-format style lingers on them.
foreach ($r in $textPriceRows) { $ws.Range("D$r").Style = "Normal" }

